$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13 ("A place to call home"): corrected narrative text + kingdom instructions, kingdom-level bump ---
$ws.Cells.Item(13, 3).Value = 'You are exhausted, emotionally, mentally and physically. You fight, you kill, you train, you craft and enchant. You sit under a tree with a small fire burning as the moon rises high in the night sky. You notice how bright, how clear, how large the moon is in the night sky.<br /> <br /> As you sit and mull over the last few months, your mind wonders back to what The Poet stated, about the abandoned house, your house, your family house.<br /> <br /> Are they dead? Are they alive? Where are they? Why did they leave the farm? So many questions, yet no answers.<br /> <br /> A slight breeze kicks up and the flames flicker and dance in the night.<br /> <br /> A voice comes from the shadows of the night, beyond the flames light and reach.<br /> <br /> “Child, I am back.” The Guide steps from the shadows.<br /> <br /> You look up from your fire and see that he is young again, keeps changing between young and old. Alas he looks bruised, broken and is limping.<br /> <br /> You ask him what happened but he shrugs it off. “It is not important. What is though is that you put your feet up. Gather some soldiers, build some walls, create a home for your self.”<br /> <br /> You tell him what The Poet said about your parents.<br /> <br /> “Silly child, there is so much for you to learn. Alas for now, Lets get you to settle down and call some place home.”<br /> <br /> You insist he explains more, but he refuse’s and starts to get angry at your demands.<br /> <br /> “Now you listen here, there is a darkness descending on this world, one escaping the pits of hell, some of it I am sure you have learned about, others you don’t know yet, but soon will. All will be answered in time. But for now, you have a job to do.”<br /> <br /> You think for a moment of refusing until you get the answers, but a swath of land you can call your own does sound enticing.'
$ws.Cells.Item(13, 4).Value = '<p>Lets take a break from questing, crafting, enchanting and fighting and lets focus on kingdoms and what they can provide for you.</p><p>You can settle a <a href="/information/kingdoms" target="_blank">kingdom</a> for free once, assuming you have no kingdoms.</p><p>Each kingdom after the first one will cost 10,000 Gold x the amount of kingdoms you already own.</p><p>Your first kingdom, again assuming you have no kingdoms, will have a 7 day protection on it, meaning no one can attack it.</p><p><strong>Note:</strong> the quest requires a kingdom level, this means you need you need to level the buildings, by resources or gold, over time to reach the required objective. You can level one building or multiple as the kingdom level is a sum of all the building levels in one kingdom.</p><p><strong>Note: The quest requires you to purchase with gold or resources units. If you purchased your buildings (2 at level 10) to will notice that you have no population. On the right hand side under Kingdom Actions for desktop or under kingdoms details for mobile, you will find a Buy Population button that lets you buy people for gold. Caution as to not exceed the max amount of population your kingdom can have or an angry NPC might come knocking.</strong></p><p><strong>Final Note: Kingdom resource replenishment and treasury interest for deposited gold occurs every hour and can be seen in the Server Message section if you go to your </strong><a href="//information/settings" target="_blank"><strong>settings</strong></a><strong> and enable the messages if they wish. <em>This can get quite spammy so it is disabled by default.</em></strong></p><p><strong>Desktop</strong></p><p><strong>- </strong>Find a place on either Surface, Labyrinth or Dungeons where you can access and what to settle a kingdom.</p><p>- Click Settle Kingdom under the map and give your kingdom a name.</p><p>- Now you can click on the Kingdoms tab and see your kingdom</p><p>- From there you can click on the kingdom name and see further details.</p><p>- On the map you will see an icon on the map that represents your kingdom and you can also click on that to get basic details as well as teleport to it, assuming you have the gold to teleport the distance.</p><p>- When you click on the kingdom via the kingdoms tab, you can see information like buildings and units as well as kingdom resources that will replenish over time, every hour.</p><p>- Lets recruit some units and upgrade our buildings which can be done with either resources, as stated, or gold. Both of which will require (real world) time to recruit the units and build the buildings.</p><p><strong>Mobile</strong></p><p><strong>- </strong>Find a place on either Surface, Labyrinth or Dungeons where you can access and what to settle a kingdom.</p><p>- In the Actions drop down, select Map Movement</p><p>- Tap Settle Kingdom under the map and give your kingdom a name.</p><p>- Now you can tap on the Kingdoms tab and see your kingdom</p><p>- From there you can tap on the kingdom name and see further details.</p><p>- On the map you will see your kingdom icon, but due to how the map works on mobile, you wont be able to tap on it for basic details unlike on desktop. How ever if you tap Teleport you can choose from a list of your kingdoms on that Plane to teleport to.</p><p>- When you tap on the kingdom via the kingdoms tab, you can see information like buildings and units as well as kingdom resources that will replenish over time, every hour.</p><p>- Lets recruit some units and upgrade our buildings which can be done with either resources, as stated, or gold. Both of which will require (real world) time to recruit the units and build the buildings.</p>'
$ws.Cells.Item(13, 27).Value = 20

# --- Row 14 ("Growing the population"): higher kingdom-level / unit requirements + new passive-level requirement ---
$ws.Cells.Item(14, 27).Value = 40
$ws.Cells.Item(14, 28).Value = 500
$ws.Cells.Item(14, 30).Value = 1

# --- Row 17: New quest 16: Discovering the truth pt. 1 ---
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 'Discovering the truth pt. 1'
$ws.Cells.Item(17, 3).Value = 'You manage to awaken from your dreams, dreams of the void. Dreams of the nothingness that haunts your mind. Haunts you in a way where you question your own existence. One that isn’t even yours, if you understood correctly.<br /> <br /> You leave the room of your INN and head down stairs to the bar area for the food. There is a message for you, left by The Guide. The Barmaid hands you the note as she brings over your food.<br /> <br /> “Child” The letter starts.<br /> <br /> “I have found answers to your questions, and solutions to our current predicament with The Creator. Alas I need you to get much, much stronger – to do this, lets trade in some books and visit a friendly goblin. By the time you are done, I shall be back.<br /> <br /> Signed,<br /> <br /> Guide.”<br /> <br /> Goblin? I have to speak to one of those foul creatures?'
$ws.Cells.Item(17, 4).Value = '<p>Lets start diving into the game now that we have some of the basics under our belt, lets visit the goblin for some quests.</p><p>Go to surface and under the Quests tab you will see - to the right, a quest called Goblins Lust for Gold. This is the story of the goblin as well as some story around The Creator for more context.</p><p>This should be an easy quest to do.</p><p>The item you get from the quest, will give additional faction points per kill, this will help you with other quests in the game where specific faction levels are a requirement.</p><p><em>Note: The Creator is not the developer of the game, this is an in game NPC</em></p>'
$ws.Cells.Item(17, 20).Value = 'Faction Hunter'
$ws.Cells.Item(17, 21).Value = 'Weaponsmith''s Master Book'
$ws.Cells.Item(17, 22).Value = 'Blacksmith''s Master Recipes'
$ws.Cells.Item(17, 41).Value = 10000
$ws.Cells.Item(17, 42).Value = 25

# --- Row 18: New quest 17: Discovering the truth pt. 2 ---
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'Discovering the truth pt. 2'
$ws.Cells.Item(18, 3).Value = 'You spoke with The Goblin and completed some tasks for some rewards. Alas you return to the INN with a new set of books and The Guide is not here.<br /> <br /> You wait well into the night, and he still doesn’t show.'
$ws.Cells.Item(18, 4).Value = '<p>Might as well head back to that goblin.</p><p>The item you get from this quest chain will help you not only hit better, much much better, but also gain levels - which will make the skill requirement easier.</p>'
$ws.Cells.Item(18, 6).Value = 'Accuracy'
$ws.Cells.Item(18, 7).Value = 100
$ws.Cells.Item(18, 20).Value = 'Goblins Accuracy'
$ws.Cells.Item(18, 41).Value = 10000
$ws.Cells.Item(18, 42).Value = 25

# --- Row 19: New quest 18: Discovering the truth pt. 3 ---
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'Discovering the truth pt. 3'
$ws.Cells.Item(19, 3).Value = 'You decided to check in on that stretch of land you call your home. Your advisors have told you about a ring in which The Goblin might be willing to trade with you.<br /> <br /> You inquire more about it and they tell you about how your ability to be a king (or queen) is shown through your kingmanship.<br /> <br /> You head out to investigate, all while handling a few bits and bobs around the kingdom.'
$ws.Cells.Item(19, 4).Value = '<p>One of your training skills is called Kingmanship, which helps in recruiting larger amounts of units and reducing time around both recruiting and upgrading as well as the cost.</p><p>The Goblin has a fun item for you to help level this skill.</p>'
$ws.Cells.Item(19, 6).Value = 'Kingmanship'
$ws.Cells.Item(19, 7).Value = 15
$ws.Cells.Item(19, 20).Value = 'Leading Armies'
$ws.Cells.Item(19, 27).Value = 80
$ws.Cells.Item(19, 28).Value = 1000
$ws.Cells.Item(19, 29).Value = 'Kingdom Management'
$ws.Cells.Item(19, 30).Value = 3
$ws.Cells.Item(19, 39).Value = 1000
$ws.Cells.Item(19, 41).Value = 50000
$ws.Cells.Item(19, 42).Value = 75

# --- Row 20: New quest 19: Discovering the truth pt. 4 ---
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'Discovering the truth pt. 4'
$ws.Cells.Item(20, 3).Value = 'You have been speaking with the local Blacksmith of the little village you are building, when The Guide appears. You excuse your self from the conversation and walk off to the side with The Guide.<br /> <br /> “You have a nice little hamlet here.” he states.<br /> <br /> You ask him, with hushed and angry tones where the hell he has been.<br /> <br /> “I have been to the land where the Shade Lord calls home. The Shadow Planes.”<br /> <br /> What could possibly be down there?<br /> <br /> “There is rumors of a cloak that makes it’s wearer invisible.”<br /> <br /> You ask him if you are to fetch it for him and he shakes his head and starts to laugh.<br /> <br /> “Child you still depend on the shop for all your shabby gear. You need much more powerful gear.” He walks off and you follow him.'
$ws.Cells.Item(20, 4).Value = '<p>I really hope you have been crafting this whole time like I told you to be long ago. This is the painful quest if you haven’t been.</p><p>The reason the levels are selected is because, you will unlock the first item – of many – that cannot be bought in the shop.</p><p>The books we got you, a couple quests back, will help a bit. If you have been enchanting, then search for weapon/armor crafting and create your first <a href="/information/equipment-sets" target="_blank">set</a> of armor crafting and weapon crafting enchantments to make it even easier and faster.</p><p><strong>Hint, I pray you have been enchanting as well child ...</strong></p>'
$ws.Cells.Item(20, 6).Value = 'Weapon Crafting'
$ws.Cells.Item(20, 7).Value = 225
$ws.Cells.Item(20, 8).Value = 'Armour Crafting'
$ws.Cells.Item(20, 9).Value = 294
$ws.Cells.Item(20, 17).Value = 'Surface'
$ws.Cells.Item(20, 18).Value = 5
$ws.Cells.Item(20, 41).Value = 2000000000
$ws.Cells.Item(20, 42).Value = 100

# --- Row 21: New quest 20: Discovering the truth pt. 5 ---
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = 'Discovering the truth pt. 5'
$ws.Cells.Item(21, 3).Value = 'The Guide has been staying in your little village for the last couple weeks, while you work day and day, night after night – crafting and creating Armour and weapons that even the local merchants don’t sell.<br /> <br /> “Now, enchant it all.” The Guide states as you two are having dinner and discussing your latest accomplishments.<br /> <br /> “I would also pay attention to your kingdom.” He states before going back to his dinner. You stare at him with rage in your eyes.'
$ws.Cells.Item(21, 4).Value = '<p>Now its time to enchant – the better the enchant on the gear, the better chance you have to survive. There are <a href="/information/enchanting?table[filters][types]=16" target="_blank">enchantments</a> that can help with enchanting.</p><p>At this stage, you should have a grasp on the core concept of the game: Gear. The Better the gear you can craft, the better enchantment and mixture of enchantments that you apply to your gear all determines if you can kill a monster in one hit.</p><p>Soon you will be going to places where creatures are much harder, where your gear matters, where your skills need to be leveled up and where you should be experimenting with <a href="/information/class-ranks" target="_blank">class ranks</a>, mixing and matching different classes abilities together and leveling them up for even more power.</p><p>Keep crafting, keep enchanting and keep leveling and soon you can take on the first challenge before you:<a href="/information/celestials" target="_blank"> Celestial Entities</a></p>'
$ws.Cells.Item(21, 6).Value = 'Enchanting'
$ws.Cells.Item(21, 7).Value = 200
$ws.Cells.Item(21, 41).Value = 2000000000
$ws.Cells.Item(21, 42).Value = 100

# --- Row 22: New quest 21: Discovering the truth pt. 6 ---
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'Discovering the truth pt. 6'
$ws.Cells.Item(22, 3).Value = 'You have been telling your advisors on how you want to expand, when a merchants comes walking into the meeting room, completely uninvited.<br /> <br /> “Who are you?!” Exclaims one of the advisors.<br /> <br /> The guards bust in and attempt to grab the man who dodges their attempts.<br /> <br /> “I am the Merchant!” he exclaims. You wave your hand and your guards stand down.<br /> <br /> Your advisors attempt to interject but you refuse to listen and instead turn and as The Merchant who and why he is here.<br /> <br /> “I have come to teach you the story of the wondrous world of Alchemy. Alas you need something first. Now where is my room?”<br /> <br /> Alchemy? The Merchant? Like The Poet?<br /> <br /> “That’s right! Like The Poet!” he states in a happy boasted voice.<br /> <br /> What is going on.'
$ws.Cells.Item(22, 4).Value = '<p>I want you to learn about Alchemy. <a href="/information/alchemy" target="_blank">Alchemy</a> is good for three things:</p><p>- Boons, which apply effects to your character for a limited amount of time but can greatly power your character.</p><p>- Kingdom Destroying Items that you can drop on other players kingdoms to deal large amounts of damage.</p><p>- Holy Oils, which you can apply to your gear later on, to give your self even more power.</p><p>You will need to complete a One Off quest on surface first to enable the new skill to appear in your craft/enchant section.</p><p>Once this is quest is done you will be able to select Alchemy and then select the item you want to craft.</p><p>Alchemy needs Gold Dust which you can get from disenchanting enchanted items. It also needs Shards, which an only be gotten from a location you cannot access yet or from fighting Celestials, which you are not strong enough for. Instead lets use what shards you do have to level Alchemy just a little bit.</p><p>Alchemy items can be seen in your inventory, under the Type Drop down called: Usable. These can be listed on the market or destroyed</p><p>We will come back to this skill later.</p>'
$ws.Cells.Item(22, 6).Value = 'Alchemy'
$ws.Cells.Item(22, 7).Value = 5
$ws.Cells.Item(22, 20).Value = 'I dream of Alchemy'
$ws.Cells.Item(22, 26).Value = 2
$ws.Cells.Item(22, 27).Value = 150
$ws.Cells.Item(22, 28).Value = 2000
$ws.Cells.Item(22, 29).Value = 'Kingdom Management'
$ws.Cells.Item(22, 30).Value = 4
$ws.Cells.Item(22, 41).Value = 1000000000
$ws.Cells.Item(22, 42).Value = 100

# --- Column width tweaks (best-fit recalculated by Excel for the longer/shorter text) ---
$ws.Columns.Item(3).ColumnWidth = 2193.3333333333335
$ws.Columns.Item(4).ColumnWidth = 4370.0
$ws.Columns.Item(21).ColumnWidth = 29.666666666666668
$ws.Columns.Item(22).ColumnWidth = 32.166666666666664
